$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet
$ws.Name = "BetaFiberA"

# 2. Fix the tiny numerical discrepancy in L13
$ws.Range("L13").Value = 0.9934699647258994

# 3. Add new row 16 of data (Gaussian Quadrature results for HKL index 14 / HexGrid-60degTilt5degRes)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.115166809739337
$ws.Range("D16").Value = 0.8915662422063955
$ws.Range("E16").Value = 1.023812253641894
$ws.Range("F16").Value = 0.9546689374693714
$ws.Range("G16").Value = 1.115166809739337
$ws.Range("H16").Value = 0.8915662422063955
$ws.Range("I16").Value = 1.041976197050596
$ws.Range("J16").Value = 0.9555254482565383
$ws.Range("K16").Value = 1.024786699429239
$ws.Range("L16").Value = 0.9128234721908246
$ws.Range("M16").Value = 1.115166809739337
$ws.Range("N16").Value = 0.9576892479241449
$ws.Range("O16").Value = 0.9963035607642496
$ws.Range("P16").Value = 0.9900407574980246

# Apply the same style as A15 (bold, bordered, centered) to A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
